$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()
try {
    $excel.Goto($ws.Range("A70"), $false)
    Write-Host "goto ok"
} catch {
    Write-Host "ERROR goto: $_"
}
$ws.Range("C94").Select()
